$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-115 is being updated from serial date
# 45172 (2023-09-03) to 45175 (2023-09-06) across the whole data range.
$ws.Range("C2:C115").Value = 45175
